$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F10").Value = 1445
$ws1.Range("F12").Value = 37
$ws1.Range("F13").Value = 366
$ws1.Range("F16").Value = 7
$ws1.Range("F21").Value = 211
$ws1.Range("F22").Value = 200

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F11").Value = 1445
$ws4.Range("F13").Value = 37
$ws4.Range("F14").Value = 366
$ws4.Range("F17").Value = 7
$ws4.Range("F22").Value = 211
$ws4.Range("F23").Value = 200
